# Add a new "Password" column (E) with per-student password values,
# each rendered as a hyperlink (mirroring the existing Gmail ID hyperlink
# column), and a header matching the style of the other column headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header cell - copy formatting from an existing header cell (A1) so it
# reuses the same bold/bordered header style.
$ws.Range("E1").Value = "Password"
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# Data cells with hyperlinks, styled like the existing Gmail ID hyperlinks
$ws.Range("E2").Value = "dhruv@2015"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:dhruv@2015")
$ws.Range("E2").Style = "Hyperlink"

$ws.Range("E3").Value = "pushp@2013"
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:pushp@2013")
$ws.Range("E3").Style = "Hyperlink"

$ws.Range("E4").Value = "git@2011"
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:git@2011")
$ws.Range("E4").Style = "Hyperlink"

# Move the active selection to mirror where the author ended up working
$ws.Range("E14").Select()
